# Applies the "upload table feature" edit to the Financial Information sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Information")

# Company name
$ws.Range("B1").Value = "TATA STEEL LIMITED"

# Net Sales (row 9)
$ws.Range("B9").Value = 1290066200000
$ws.Range("C9").Value = 1290213500000

# Other Income (row 11)
$ws.Range("B11").Value = 33254800000
$ws.Range("C11").Value = 14520200000

# Total Income (row 12)
$ws.Range("B12").Value = 1323321000000
$ws.Range("C12").Value = 1304733700000

# Dividend received (row 19)
$ws.Range("B19").Value = 2439200000
$ws.Range("C19").Value = 2439200000

# Interest received (row 20)
$ws.Range("B20").Value = 27207100000
$ws.Range("C20").Value = 9430000000

# Interest paid (row 34)
$ws.Range("B34").Value = 37921400000
$ws.Range("C34").Value = 27920800000

# Employee benefit expense (row 66)
$ws.Range("B66").Value = 66162900000
$ws.Range("C66").Value = 63658000000
